$wb = $excel.ActiveWorkbook

# This script applies a batch of market-price refresh values to the
# H:N "current price / profit" columns across several Leve sheets.
# Values were recomputed from an updated market-data snapshot; most
# cells are simple value replacements. A couple of rows on the CUL
# sheet (81, 84) also drop a now-unused LeveProfitNQ (M) cell, shifting
# LevePriceHQ (L) to absorb the new combined value.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6224.5
$ws.Range("J40").Value = 4999.5
$ws.Range("L40").Value = 4999.5
$ws.Range("N40").Value = -5349.5

$ws.Range("H62").Value = 3835.3333
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 4003
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 4003
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -5251

$ws.Range("H65").Value = 3835.3333
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 4003
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 20015
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -26255

$ws.Range("H98").Value = 4660.8
$ws.Range("I98").Value = 4660.8
$ws.Range("K98").Value = 4660.8
$ws.Range("M98").Value = -3162.8

$ws.Range("H112").Value = 2116.6453
$ws.Range("J112").Value = 2116.6453
$ws.Range("L112").Value = 6349.9359
$ws.Range("N112").Value = -8565.9359

$ws.Range("H113").Value = 3879.4
$ws.Range("J113").Value = 3999.3333
$ws.Range("L113").Value = 3999.3333
$ws.Range("N113").Value = -10507.3333

$ws.Range("H116").Value = 9979.166999999999
$ws.Range("J116").Value = 9979.166999999999
$ws.Range("L116").Value = 9979.166999999999
$ws.Range("N116").Value = -16863.167

$ws.Range("H122").Value = 4660.8
$ws.Range("I122").Value = 4660.8
$ws.Range("K122").Value = 13982.4
$ws.Range("M122").Value = -11532.4

$ws.Range("H138").Value = 319006.4
$ws.Range("I138").Value = 2320.9395
$ws.Range("J138").Value = 546193.8
$ws.Range("K138").Value = 6962.818499999999
$ws.Range("L138").Value = 1638581.4
$ws.Range("M138").Value = -1822.818499999999
$ws.Range("N138").Value = -1648861.4


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 836.0952
$ws.Range("I2").Value = 550.38464
$ws.Range("J2").Value = 1300.375
$ws.Range("K2").Value = 550.38464
$ws.Range("L2").Value = 1300.375
$ws.Range("M2").Value = -437.38464
$ws.Range("N2").Value = -1526.375

$ws.Range("H32").Value = 4673.9536
$ws.Range("I32").Value = 3570.4062
$ws.Range("J32").Value = 7884.273
$ws.Range("K32").Value = 3570.4062
$ws.Range("L32").Value = 7884.273
$ws.Range("M32").Value = -3283.4062
$ws.Range("N32").Value = -8458.273000000001

$ws.Range("H45").Value = 36793.418
$ws.Range("I45").Value = 39964.273
$ws.Range("K45").Value = 39964.273
$ws.Range("M45").Value = -39587.273

$ws.Range("H63").Value = 4173.7334
$ws.Range("I63").Value = 4373.615
$ws.Range("J63").Value = 2874.5
$ws.Range("K63").Value = 4373.615
$ws.Range("L63").Value = 2874.5
$ws.Range("M63").Value = -3687.615
$ws.Range("N63").Value = -4246.5

$ws.Range("H66").Value = 4173.7334
$ws.Range("I66").Value = 4373.615
$ws.Range("J66").Value = 2874.5
$ws.Range("K66").Value = 21868.075
$ws.Range("L66").Value = 14372.5
$ws.Range("M66").Value = -18436.075
$ws.Range("N66").Value = -21236.5

$ws.Range("H74").Value = 294225.4
$ws.Range("I74").Value = 348681.3
$ws.Range("K74").Value = 348681.3
$ws.Range("M74").Value = -347807.3

$ws.Range("H77").Value = 294225.4
$ws.Range("I77").Value = 348681.3
$ws.Range("K77").Value = 1743406.5
$ws.Range("M77").Value = -1739038.5

$ws.Range("H116").Value = 836.0952
$ws.Range("I116").Value = 550.38464
$ws.Range("J116").Value = 1300.375
$ws.Range("K116").Value = 550.38464
$ws.Range("L116").Value = 1300.375
$ws.Range("M116").Value = 1743.61536
$ws.Range("N116").Value = -5888.375

$ws.Range("H122").Value = 3387.9666
$ws.Range("I122").Value = 3017.24
$ws.Range("K122").Value = 9051.719999999999
$ws.Range("M122").Value = -6601.719999999999

$ws.Range("H132").Value = 2043.5454
$ws.Range("I132").Value = 1098.32
$ws.Range("J132").Value = 4997.375
$ws.Range("K132").Value = 3294.96
$ws.Range("L132").Value = 14992.125
$ws.Range("M132").Value = -764.96
$ws.Range("N132").Value = -20052.125

$ws.Range("H134").Value = 101995.664
$ws.Range("J134").Value = 101995.664
$ws.Range("L134").Value = 101995.664
$ws.Range("N134").Value = -112135.664


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 836.0952
$ws.Range("I3").Value = 550.38464
$ws.Range("J3").Value = 1300.375
$ws.Range("K3").Value = 550.38464
$ws.Range("L3").Value = 1300.375
$ws.Range("M3").Value = -436.38464
$ws.Range("N3").Value = -1528.375

$ws.Range("H107").Value = 1673.6666
$ws.Range("I107").Value = 1414.909
$ws.Range("K107").Value = 1414.909
$ws.Range("M107").Value = 505.0909999999999

$ws.Range("H134").Value = 4122.2964
$ws.Range("I134").Value = 3991.279
$ws.Range("J134").Value = 4634.4546
$ws.Range("K134").Value = 11973.837
$ws.Range("L134").Value = 13903.3638
$ws.Range("M134").Value = -9438.837
$ws.Range("N134").Value = -18973.3638


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3674.279
$ws.Range("I31").Value = 2265.5588
$ws.Range("K31").Value = 2265.5588
$ws.Range("M31").Value = -1970.5588

$ws.Range("H34").Value = 3674.279
$ws.Range("I34").Value = 2265.5588
$ws.Range("K34").Value = 2265.5588
$ws.Range("M34").Value = -2063.5588

$ws.Range("H58").Value = 2795.4138
$ws.Range("I58").Value = 2182.5715
$ws.Range("J58").Value = 3367.4
$ws.Range("K58").Value = 2182.5715
$ws.Range("L58").Value = 3367.4
$ws.Range("M58").Value = -1979.5715
$ws.Range("N58").Value = -3773.4

$ws.Range("H99").Value = 6299.3335
$ws.Range("I99").Value = 4999.5
$ws.Range("J99").Value = 6949.25
$ws.Range("K99").Value = 4999.5
$ws.Range("L99").Value = 6949.25
$ws.Range("M99").Value = -3501.5
$ws.Range("N99").Value = -9945.25

$ws.Range("H107").Value = 617.6818
$ws.Range("I107").Value = 370.7857
$ws.Range("K107").Value = 370.7857
$ws.Range("M107").Value = 1549.2143

$ws.Range("H122").Value = 3406.2104
$ws.Range("I122").Value = 2844.7693
$ws.Range("K122").Value = 8534.3079
$ws.Range("M122").Value = -6084.3079

$ws.Range("H126").Value = 6299.3335
$ws.Range("I126").Value = 4999.5
$ws.Range("J126").Value = 6949.25
$ws.Range("K126").Value = 14998.5
$ws.Range("L126").Value = 20847.75
$ws.Range("M126").Value = -12528.5
$ws.Range("N126").Value = -25787.75

$ws.Range("H132").Value = 2858.606
$ws.Range("I132").Value = 2388.0952
$ws.Range("J132").Value = 3682
$ws.Range("K132").Value = 7164.285600000001
$ws.Range("L132").Value = 11046
$ws.Range("M132").Value = -4634.285600000001
$ws.Range("N132").Value = -16106

$ws.Range("H136").Value = 2795.4138
$ws.Range("I136").Value = 2182.5715
$ws.Range("J136").Value = 3367.4
$ws.Range("K136").Value = 6547.7145
$ws.Range("L136").Value = 10102.2
$ws.Range("M136").Value = -3997.7145
$ws.Range("N136").Value = -15202.2


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5417
$ws.Range("I3").Value = 4047.5557
$ws.Range("K3").Value = 12142.6671
$ws.Range("M3").Value = -12030.6671

$ws.Range("H60").Value = 1114750.1
$ws.Range("I60").Value = 3334853.8
$ws.Range("K60").Value = 10004561.4
$ws.Range("M60").Value = -10004310.4

$ws.Range("H81").Value = 4867.875
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 4867.875
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 14603.625
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -16849.625

$ws.Range("H84").Value = 4867.875
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 4867.875
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 43810.875
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -55042.875


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1336.6923
$ws.Range("I102").Value = 786.44446
$ws.Range("K102").Value = 786.44446
$ws.Range("M102").Value = 835.55554


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4149
$ws.Range("I7").Value = 2299
$ws.Range("K7").Value = 2299
$ws.Range("M7").Value = -2187

$ws.Range("H22").Value = 1000

$ws.Range("H27").Value = 1000

$ws.Range("H40").Value = 257493.25
$ws.Range("I40").Value = 999999
$ws.Range("J40").Value = 9991.333000000001
$ws.Range("K40").Value = 999999
$ws.Range("L40").Value = 9991.333000000001
$ws.Range("M40").Value = -999863
$ws.Range("N40").Value = -10263.333

$ws.Range("H55").Value = 313.83334
$ws.Range("I55").Value = 154.5
$ws.Range("K55").Value = 154.5
$ws.Range("M55").Value = 18.5

$ws.Range("H95").Value = 34998.5
$ws.Range("J95").Value = 34998.5
$ws.Range("L95").Value = 34998.5
$ws.Range("N95").Value = -40490.5

$ws.Range("H126").Value = 4149
$ws.Range("I126").Value = 2299
$ws.Range("K126").Value = 6897
$ws.Range("M126").Value = -4427


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 111118580
$ws.Range("I136").Value = 142859310
$ws.Range("K136").Value = 428577930
$ws.Range("M136").Value = -428575380

